# Update "想去人数" (number of people wanting to go) counts that changed
# between the two generated-data snapshots.
# 展览(F2): 152 -> 153
# 展览(F3): 43  -> 44
# 全部类型(F2): 152 -> 153
# 全部类型(F3): 43  -> 44

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 153
    $ws.Range("F3").Value = 44
}
